$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: "*_old" -> "*_FV2410", "*_new" -> "*_FV2504" ---
$oldHeaders = @("A1","B1","C1","D1","E1","F1","G1","H1","I1","J1")
foreach ($addr in $oldHeaders) {
    $cell = $ws.Range($addr)
    $cell.Value2 = ($cell.Value2 -replace "_old$", "_FV2410")
}

$newHeaders = @("L1","M1","N1","O1","P1","Q1","R1","S1","T1","U1")
foreach ($addr in $newHeaders) {
    $cell = $ws.Range($addr)
    $cell.Value2 = ($cell.Value2 -replace "_new$", "_FV2504")
}

# --- 2. Turn the used range into an Excel Table (adds xl/tables/table1.xml + tableParts) ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U89"), 0, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ---
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
